# SC2x W1L2 NERD - re-solve the NERD3 scenario so a single DC (WO) covers
# all demand (min DC = 1, max DC = 1) instead of the previous 4-DC solution.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NERD3")
$ws.Activate()

# --- Solver bookkeeping: "max number of DC" constraint dropped from 5 to 1,
#     and the solver_num named range (NERD3) recounts from 8 to 6 constraints.
$ws.Range("E25").Value = 1

$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "NERD3!solver_num") {
        $n.RefersTo = "=6"
    }
}

# --- Decision variables (C17:C21): only the 5th candidate (WO, row 21) is
#     selected now; all others are deselected.
$ws.Range("C17").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("C21").Value = 1

# --- Shipment matrix (C6:N10): everything that used to be spread across
#     rows 6-9 now ships entirely out of row 10 (WO), so zero out rows 6-9
#     and move the totals onto row 10.
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N")

foreach ($col in $cols) {
    $ws.Range($col + "6").Value = 0
    $ws.Range($col + "7").Value = 0
    $ws.Range($col + "8").Value = 0
    $ws.Range($col + "9").Value = 0
}

$row10 = @{
    "C" = 425.00000000000011
    "D" = 12
    "E" = 43
    "F" = 125
    "G" = 110
    "H" = 86.000000000000057
    "I" = 129
    "J" = 28
    "K" = 66
    "L" = 320.00000000000006
    "M" = 220
    "N" = 181.99999999999997
}

foreach ($col in $cols) {
    $ws.Range($col + "10").Value = $row10[$col]
}

# --- Restore the selection to B1, matching where the user left off.
$ws.Range("B1").Select()
